# BAU Fraction of CCS Potential Achieved by Industry.xlsx
# - Rename the "NoSettings" setting-suffix used throughout the BAU Emissions
#   label column to "test".
# - Refresh the BAU Emissions[natural gas if,iron and steel] energy series
#   (row 94) with updated model output.
# - Bump the "last updated" date on the About sheet.
# - Leave the workbook positioned on the About tab (previously the
#   "Current and Planned Capacity" tab was active) and restore the BAU
#   Emissions sheet's own scroll/selection state.

$wb = $excel.ActiveWorkbook

# --- BAU Emissions sheet: label text + refreshed data -----------------
$wsBau = $wb.Worksheets.Item("BAU Emissions")
$wsBau.Activate()

# Every label in column A that ends " : NoSettings" becomes " : test".
[void]$wsBau.Range("A1:A300").Replace(" : NoSettings", " : test")

# Row 94 = "Industrial Sector Energy Related Emissions before CCS[natural
# gas if,iron and steel 241,CO2]" — updated 2032-2050 series.
$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300
$wsBau.Range("Q94:AE94").Value = 5005380

# Restore this sheet's own selection/scroll state.
[void]$wsBau.Range("A30:AE280").Select()

# --- About sheet: bump the last-updated date, make it the active tab --
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387
$wsAbout.Activate()
